$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new contingency lines ("line7"/"line8") are inserted right after the
# existing "line1".."line6" rows. That pushes the whole "extr1".."extr8"
# block down by two rows (old row 8 -> new row 10, ... old row 15 -> new
# row 17), so every extr row's id (col A) and label (col B) shift by two
# positions as well. New rows 16/17 need the same bold/centered/bordered
# style as the rest of the A-column id cells - copy it down from A15.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null

# (row, A, B, C, D, E) for every row from "line7" through the new "extr8".
$rows = @(
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
